$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B10: was stored as text "3", should be a real number 3
$ws.Range("B10").Value = 3

# Add new row 11 with payment record data.
# Some of these columns (member_id, the "date" column) must stay TEXT
# rather than being auto-converted to number/date by COM, so force the
# cell format to Text first (mirrors typing an apostrophe-prefixed value).
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "4"
$ws.Range("D11").Value = 10000
$ws.Range("G11").Value = "Paid"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "2025-03-17"
$ws.Range("I11").Value = "Amna"
$ws.Range("J11").Value = "Super"

# A11 / C11 / E11 / F11 stay blank (matching the blank id / package_id /
# payment_date / next_due_date columns used throughout the sheet).
# Touching NumberFormat materialises the otherwise-empty cell without
# coercing it to a typed value.
"A11", "C11", "E11", "F11" | ForEach-Object {
    $ws.Range($_).NumberFormat = "general"
}
